$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make sure new rows (11, 12) carry the same style as the existing rows ---
# (column A uses a bold/centered/thin-bordered style; copy that formatting
#  from the last existing row onto the two new rows)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)

# --- Set A, B, D, F, G columns for all data rows (2-12) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("D2").Value = 239.5
$ws.Range("F2").Value = 114.0541439806413
$ws.Range("G2").Value = 5
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("D3").Value = 224.5
$ws.Range("F3").Value = 112.1024382991377
$ws.Range("G3").Value = 14.5
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("D4").Value = 225.5
$ws.Range("F4").Value = 112.1172077922078
$ws.Range("G4").Value = 14.5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("D5").Value = 238
$ws.Range("F5").Value = 116
$ws.Range("G5").Value = 1
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("D6").Value = 222
$ws.Range("F6").Value = 114.2004310344828
$ws.Range("G6").Value = 2.5
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2023
$ws.Range("D7").Value = 220.5
$ws.Range("F7").Value = 112.780701754386
$ws.Range("G7").Value = 1.5
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2023
$ws.Range("D8").Value = 226
$ws.Range("F8").Value = 112.3834586466165
$ws.Range("G8").Value = 5
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 2023
$ws.Range("D9").Value = 237
$ws.Range("F9").Value = 116.6393939393939
$ws.Range("G9").Value = 3
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 2023
$ws.Range("D10").Value = 232.5
$ws.Range("F10").Value = 114.3299240210403
$ws.Range("G10").Value = 7.5
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 2023
$ws.Range("D11").Value = 233.5
$ws.Range("F11").Value = 115.7069805194805
$ws.Range("G11").Value = 4
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 2023
$ws.Range("D12").Value = 238
$ws.Range("F12").Value = 115.9553571428571
$ws.Range("G12").Value = 2.5

# --- Set H column (home_team) for all rows, then I column (away_team) for all rows ---
# (column-major order matters for shared-string placement parity with source)
$ws.Range("H2").Value = "Charlotte"
$ws.Range("H3").Value = "Cleveland"
$ws.Range("H4").Value = "Philadelphia"
$ws.Range("H5").Value = "Indiana"
$ws.Range("H6").Value = "NewYork"
$ws.Range("H7").Value = "Miami"
$ws.Range("H8").Value = "Chicago"
$ws.Range("H9").Value = "OklahomaCity"
$ws.Range("H10").Value = "Dallas"
$ws.Range("H11").Value = "GoldenState"
$ws.Range("H12").Value = "Portland"
$ws.Range("I2").Value = "Atlanta"
$ws.Range("I3").Value = "SanAntonio"
$ws.Range("I4").Value = "Houston"
$ws.Range("I5").Value = "Utah"
$ws.Range("I6").Value = "Brooklyn"
$ws.Range("I7").Value = "Denver"
$ws.Range("I8").Value = "Orlando"
$ws.Range("I9").Value = "NewOrleans"
$ws.Range("I10").Value = "Minnesota"
$ws.Range("I11").Value = "Washington"
$ws.Range("I12").Value = "LALakers"

# --- Set remaining numeric columns J..AA for all rows ---
$ws.Range("J2").Value = 0.5220810647307925
$ws.Range("K2").Value = 100.1416212946158
$ws.Range("L2").Value = 112.9035995160315
$ws.Range("M2").Value = 116.2827888687235
$ws.Range("N2").Value = 75.55184513006657
$ws.Range("O2").Value = 0.3458312159709618
$ws.Range("P2").Value = 0.5602218693284937
$ws.Range("Q2").Value = 0.2493947368421053
$ws.Range("R2").Value = 11.14110707803993
$ws.Range("S2").Value = 12.36382335148215
$ws.Range("T2").Value = 0.2053209316394434
$ws.Range("U2").Value = 0.9978490287020233
$ws.Range("V2").Value = 0.9831696054197292
$ws.Range("W2").Value = 10.67932910905979
$ws.Range("X2").Value = 0.3836963097398669
$ws.Range("Y2").Value = 40.5
$ws.Range("Z2").Value = 76.19999999999999
$ws.Range("AA2").Value = 0.5067845343783219
$ws.Range("J3").Value = 0.5112480739599383
$ws.Range("K3").Value = 97.64485578352659
$ws.Range("L3").Value = 113.9010258697592
$ws.Range("M3").Value = 115.8981118049361
$ws.Range("N3").Value = 76.7065120428189
$ws.Range("O3").Value = 0.3562408563782337
$ws.Range("P3").Value = 0.5760346416889681
$ws.Range("Q3").Value = 0.2541992268807612
$ws.Range("R3").Value = 12.57946773713946
$ws.Range("S3").Value = 12.47756467439786
$ws.Range("T3").Value = 0.1984084151055605
$ws.Range("U3").Value = 0.9807737384001547
$ws.Range("V3").Value = 0.985679892568631
$ws.Range("W3").Value = 11.11398513741878
$ws.Range("X3").Value = 0.4363663395777579
$ws.Range("Y3").Value = 34.5
$ws.Range("Z3").Value = 75.2
$ws.Range("AA3").Value = 0.5153577947335778
$ws.Range("J4").Value = 0.5409090909090909
$ws.Range("K4").Value = 97.73595779220778
$ws.Range("L4").Value = 114.0185714285714
$ws.Range("M4").Value = 116.1429220779221
$ws.Range("N4").Value = 76.60681818181817
$ws.Range("O4").Value = 0.3854394480519481
$ws.Range("P4").Value = 0.5772524350649351
$ws.Range("Q4").Value = 0.2943868506493506
$ws.Range("R4").Value = 12.99003246753247
$ws.Range("S4").Value = 12.17079545454545
$ws.Range("T4").Value = 0.2249333603896103
$ws.Range("U4").Value = 0.9809029553124041
$ws.Range("V4").Value = 0.9850454405440763
$ws.Range("W4").Value = 11.51438712676906
$ws.Range("X4").Value = 0.4433441558441559
$ws.Range("Y4").Value = 37
$ws.Range("Z4").Value = 74.95
$ws.Range("AA4").Value = 0.5117083214312659
$ws.Range("J5").Value = 0.5263157894736842
$ws.Range("K5").Value = 99.71637931034481
$ws.Range("L5").Value = 116.0163793103448
$ws.Range("M5").Value = 117.2439655172414
$ws.Range("N5").Value = 73.47931034482758
$ws.Range("O5").Value = 0.4332155172413792
$ws.Range("P5").Value = 0.5826810344827587
$ws.Range("Q5").Value = 0.2682672413793104
$ws.Range("R5").Value = 12.62413793103448
$ws.Range("S5").Value = 12.11034482758621
$ws.Range("T5").Value = 0.2178146551724137
$ws.Range("U5").Value = 1.014873140857393
$ws.Range("V5").Value = 0.9880617998217018
$ws.Range("W5").Value = 10.38247010106923
$ws.Range("X5").Value = 0.4568965517241379
$ws.Range("Y5").Value = 24
$ws.Range("Z5").Value = 75.09999999999999
$ws.Range("AA5").Value = 0.4844894316827802
$ws.Range("J6").Value = 0.4951298701298701
$ws.Range("K6").Value = 96.79128694581283
$ws.Range("L6").Value = 117.0265086206897
$ws.Range("M6").Value = 115.1648706896551
$ws.Range("N6").Value = 75.02429187192116
$ws.Range("O6").Value = 0.385719827586207
$ws.Range("P6").Value = 0.5887358374384237
$ws.Range("Q6").Value = 0.2689963054187192
$ws.Range("R6").Value = 11.56887315270936
$ws.Range("S6").Value = 11.20668103448276
$ws.Range("T6").Value = 0.2147085899014778
$ws.Range("U6").Value = 0.999128880441669
$ws.Range("V6").Value = 0.9661284078542733
$ws.Range("W6").Value = 11.1707250111989
$ws.Range("X6").Value = 0.561884236453202
$ws.Range("Y6").Value = 44.5
$ws.Range("Z6").Value = 76.15000000000001
$ws.Range("AA6").Value = 0.4799477642716105
$ws.Range("J7").Value = 0.4649122807017544
$ws.Range("K7").Value = 97.13859649122806
$ws.Range("L7").Value = 115.5166666666667
$ws.Range("M7").Value = 113.2535087719298
$ws.Range("N7").Value = 77.75614035087722
$ws.Range("O7").Value = 0.3813684210526316
$ws.Range("P7").Value = 0.5867280701754385
$ws.Range("Q7").Value = 0.2645087719298245
$ws.Range("R7").Value = 12.30701754385965
$ws.Range("S7").Value = 13.32280701754386
$ws.Range("T7").Value = 0.2073377192982457
$ws.Range("U7").Value = 0.9867078018756429
$ws.Range("V7").Value = 1.016921609538263
$ws.Range("W7").Value = 10.0408138261072
$ws.Range("X7").Value = 0.6228070175438596
$ws.Range("Y7").Value = 50.5
$ws.Range("Z7").Value = 75.65000000000001
$ws.Range("AA7").Value = 0.48143797011866
$ws.Range("J8").Value = 0.4597402597402597
$ws.Range("K8").Value = 98.61116854636592
$ws.Range("L8").Value = 113.1957393483709
$ws.Range("M8").Value = 114.558395989975
$ws.Range("N8").Value = 77.75393170426065
$ws.Range("O8").Value = 0.3464617794486216
$ws.Range("P8").Value = 0.582845394736842
$ws.Range("Q8").Value = 0.2814584899749373
$ws.Range("R8").Value = 12.58201754385965
$ws.Range("S8").Value = 12.38618421052632
$ws.Range("T8").Value = 0.2161302474937343
$ws.Range("U8").Value = 0.9832323591130055
$ws.Range("V8").Value = 0.8890318496114061
$ws.Range("W8").Value = 11.07386678462388
$ws.Range("X8").Value = 0.4338972431077694
$ws.Range("Y8").Value = 34
$ws.Range("Z8").Value = 75.5
$ws.Range("AA8").Value = 0.4609081411531402
$ws.Range("J9").Value = 0.5807017543859649
$ws.Range("K9").Value = 99.79044657097288
$ws.Range("L9").Value = 115.6638437001595
$ws.Range("M9").Value = 114.2801435406699
$ws.Range("N9").Value = 74.95846889952153
$ws.Range("O9").Value = 0.3527744816586923
$ws.Range("P9").Value = 0.5784594896331737
$ws.Range("Q9").Value = 0.2682054226475279
$ws.Range("R9").Value = 11.87553429027113
$ws.Range("S9").Value = 13.41483253588517
$ws.Range("T9").Value = 0.2163292663476874
$ws.Range("U9").Value = 1.020467138577375
$ws.Range("V9").Value = 1.062822918699412
$ws.Range("W9").Value = 11.55173267041994
$ws.Range("X9").Value = 0.4998405103668262
$ws.Range("Y9").Value = 34.5
$ws.Range("Z9").Value = 75.90000000000001
$ws.Range("AA9").Value = 0.4938294150887351
$ws.Range("J10").Value = 0.521624780829924
$ws.Range("K10").Value = 97.86452367036821
$ws.Range("L10").Value = 115.801651081239
$ws.Range("M10").Value = 115.395601987142
$ws.Range("N10").Value = 75.21142606662769
$ws.Range("O10").Value = 0.4371237580362362
$ws.Range("P10").Value = 0.5971101694915255
$ws.Range("Q10").Value = 0.2953794564582115
$ws.Range("R10").Value = 12.32312974868498
$ws.Range("S10").Value = 12.5134862653419
$ws.Range("T10").Value = 0.2273659409701929
$ws.Range("U10").Value = 1.000261802458795
$ws.Range("V10").Value = 1.052631433805718
$ws.Range("W10").Value = 11.12350221232233
$ws.Range("X10").Value = 0.521478667445938
$ws.Range("Y10").Value = 49
$ws.Range("Z10").Value = 75.5
$ws.Range("AA10").Value = 0.4947953396634655
$ws.Range("J11").Value = 0.5508385744234801
$ws.Range("K11").Value = 99.96920454545457
$ws.Range("L11").Value = 115.1288149350649
$ws.Range("M11").Value = 115.1603571428572
$ws.Range("N11").Value = 76.51862012987013
$ws.Range("O11").Value = 0.422649025974026
$ws.Range("P11").Value = 0.5916785714285713
$ws.Range("Q11").Value = 0.248712012987013
$ws.Range("R11").Value = 12.95881493506493
$ws.Range("S11").Value = 11.31508116883117
$ws.Range("T11").Value = 0.2055767857142857
$ws.Range("U11").Value = 1.012309540852848
$ws.Range("V11").Value = 1.010713265385966
$ws.Range("W11").Value = 10.99952910247054
$ws.Range("X11").Value = 0.4863636363636363
$ws.Range("Y11").Value = 44
$ws.Range("Z11").Value = 76.19999999999999
$ws.Range("AA11").Value = 0.4755091291494521
$ws.Range("J12").Value = 0.5038986354775828
$ws.Range("K12").Value = 99.18088972431075
$ws.Range("L12").Value = 116.1247180451128
$ws.Range("M12").Value = 116.8157581453634
$ws.Range("N12").Value = 76.11126253132831
$ws.Range("O12").Value = 0.3758009085213033
$ws.Range("P12").Value = 0.590717731829574
$ws.Range("Q12").Value = 0.2975758145363409
$ws.Range("R12").Value = 12.25653195488722
$ws.Range("S12").Value = 11.04620927318296
$ws.Range("T12").Value = 0.2183666979949875
$ws.Range("U12").Value = 1.014482564679415
$ws.Range("V12").Value = 1.016498316498317
$ws.Range("W12").Value = 11.91764636535575
$ws.Range("X12").Value = 0.469141604010025
$ws.Range("Y12").Value = 42
$ws.Range("Z12").Value = 75.05
$ws.Range("AA12").Value = 0.4663094471054639
